$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename sheet from "My Series" to "Data"
$ws.Name = "Data"

# 2. Update cell A11 text from "Function Description" to "Function Information"
$ws.Range("A11").Value = "Function Information"

# 3. Change the numeric format of the data column (B27:B43) from "0.000" to "###0.000"
$ws.Range("B27:B43").NumberFormat = "###0.000"

# 4. Update the cached metadata blob stored in the A1 cell comment
$ws.Range("A1").Comment.Text("ER0AAB+LCAAAAAAAAAOlWVtvI0kV/istP4GE3d1OZjaJKr3yLVkLO45sh0z2BZW7K3GTcrfpqk7it0UCLVoWIYRm0XIVT4uQGEawKy0zXP7LapIZnvgLnLr0zXZ23GEUTbrO+U7VqVPnVhX07s2MGlckYn4Y7FfsmlUxSOCGnh9c7Fdifl61H1fedVDnxiX0GEd4RjiADZAK2N4N8/crU87ne6Z5fX1du96qhdGFWbcs23zS743cKZnhqh8wjgOXVFIp7+1SFQe1vFmfcOxhjpXkfqU76tZaxHfbQOvjAF+QqNaMmR8QxjoB97lPmJCMCOak1e5/R23Mqdce12xkrtAzZDP2qadwBaSiaxwsS8b+jDh1y96pWjvVLXts7e5Z9p5Vr+3s1t9PBFMg6mHGRyS68l1JGHE8m0txa2fLtnYt26ojcy0I5soM4KAB9YbkymfEaxFKWSmLmPoAGy6HXZczpoXMnKye6OEqHEZ4Ph37nJJyagz7TWMWaF2ySRx0EEbEBfs9SKUjcj2ItFnH8x5wx1M/4os2XpSe64SRaDAXRion6qB2GPAGJRE/mcNZEw9cARgOj2KCzHuYmVDbZy58+0FMPOccU5YXKjDRaRhdsjl2yRHEsSnmuA5oiD1wOO4z7rvZoisMdByFc5gRFm+G1DuAWTV4DSOduRuAicWyzTC8zLRbx0TyVOX5wpnOME/gK3Q0mobXg4AuRvGEuZE/IV67maDX8pAISC3dihkPZ6BFRkKKlqMs4B8E4DIZtYnrzzA9pmBE5mzBLAUCasQ8PPd5K6TxLGCJTktUdAo7GpObdIfpGA3gcANh9DDoBglemXktqygwDK/TNVcZ0gg5coO5yXGvMpbBbaAlx7fKkScidnngU6gP+bPIUYteMZoSwte6hOIgkQoPRMVxGpQiMxsicEtwbVDDsaBoVOXP2LL25A8sm7JRJ/DuxyVMdBTPBhMI3yu5IccG3hIJwRZok+LgEqinPp8eNRLV13CQ2vC9+FUegkCdU7yQ5NQoeRrqBi6NPaLivxucS48UuqkzvJeNVkg9CGkH4WAxXswhDTN/j8PHfgUK8x7jEZT+iuOGccCjhUgUyNTQt8mweBLIBTDdWOY8It+PoeNYHMSB2wq9zVfzlHVOAp9vrmEYRyr7bS4irScSYczaRKQUmeM3lnfL7IlFpeCzgMzCwHc3tzYYWWjvPWAjLImqjSWIiq+N8RTKuKpyItY3FougXYS6VmqZBmOh60tn1eHh5eTNe0KmTc5xTKFV41BRL9JUu0xGDXa5jMmT0ElEk4TniEaYQSfserOaC72C6PZqbjgTBBMa0NMRMvN40fC4pBNc9HBwEUNLkeaVZXqabkU5HEc4YGI7aQexlHnXg1CSp1Rn46wmO2QuQdCYzOZhhGkfbOIfaI/TnRE0HX3Mp3oEVYwSN7GvmYmmUkWlEp3fBpPlSO1AxLrOkEtECRL+qdrtDJPRkNhiHyKStjD1J5FKqEnRXseDs8rawCT1is2VbAkT88OVC+rst8lC9OHZQNOlt9oJQ/muyKHOaLi9U39kbcHdQo6R3PGQYGp0II45MbrBFWF8BmJ7xpAw34MvH9M94z0yIT7UP2kiXYFKS+fl0EGS4qUqDdC3SCkCoLO48KGCrAJTTibgnBEc0UUOqLbaC13A3f3437e/ef7qxWd3Hz1988UP//uPX736589un/0IPu7++rfbj3+ptqnAaIwnlEiFxs2dHWtrG/wsJSFhXFM2wV7sckk7O5O9bzpG+gonB61Ot3XYa8pUkhITcVVNTHE7XIRxNhypTciF5JGaiScoiDNOUpMeF7i56uSI29oVKaLz/PsElS1ev/zs9cs/3yutDZa1Wfbu7qOqXX9rFwb3XnsFl3ZhvUL6F+DtqvWoWq/nwEsYNFS5P7VT13Pktbq+ZdlpGvdSR14HWmbpmcb4wlySU6SW6oxSF8iPE6Z0/DGESMpWoZAbaBf94idv/vK0gNLW1ZTiLKCc7GDEYmYykFMfDcfGaHAybHWMcWck/CTj5XBq8q8B69XTeCo4VRDEmH7LgLoOdcyowN2nYoTnBsHu1FhAJObisOBs66hqoQdOuazlYRTGc3UiOYGMugaZZpO1EmtyjeRJe64knYy1Bq50vf375+sE9EbaWS+bvnTkaajAUaQcX0ftp/969eWHr168uHv+89svf1CYQa+TPgCAn0M05Yep20PK0/VmiYJOR9KYl9Z3c/VFE8Ud6jj0A84c+7G8PukRAlFbzCZ/o+4MSp6cWNoL6EsU9B5mnRuuA9s5QmaRAHrOMVTbMLtlpgSVwzO7/ue3v7v79ed3nzx/8+Gfbj/64+3Hn7x++fs3z/6gou7u6fO7nz7TWX65EEhdxN1V9X+GfAlxDRGNhqjdxlcf/MIIQm5Ay2HEMiN99cGnucmEorI5yWaGbi5VpKjCCjQvLOSMnCqpDgW5VEQ1AC1RwrZShC5i4dx3s0Xer4qpRNxJxje642rMiBFCN/VN2EkRnAlvKqdFVEk9fseq23XNVdqILUwwy5n+kIYTaDIShnxqWIIUpL5eIMPK9Q57g2ajl0GUEoPII5FwQ/WBkpZSlJQuS0aJq+UowIXGz42peB1aga2y0plzaczUzyznDU+kv/WvFAUEasVRpBqiQL/aj+I5NMPJY9z9fPlAmet/j1Svmu+Is3G3XeTDOMeFQlhkC4Lky9SkWSpNdZl40VHt7JEwTTYEXuFRE8yhH+ZVp3UFfWVkirzTiaIwWpt8Mk4C60MnDRnFzCyeYuSZqq7by84qISQJL/1Qlz69w7BNKOHlXq3NTLofXj1YFs6+rGiXDainjVnu6pGaJZsg/3QvHOX/fblXztaIImisxFNf6af25M46hKtuSW3UVqSguAHC6vot/MCPGH8iMoH+UpSzlHKmOtQn4sKlPuT4zNnWf8EBgJmf3SyomYQuV39ACWnPn/klr4VWEt/FScCW87lq4brlPEWUliNyAw1mbgZIipPvQdlQTyhlZlMOC7k0lRfPlsy/mPKyir0zwcQjE6vqTki9uu1ZO9VdQraqtg3/Y7det6xH4tFTTw6ZwyfXJRcxkwPL/rDp/A+3LVrkER0AAA==")
